# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1252
$ws1.Range("F3").Value = 17108
$ws1.Range("F5").Value = 1665
$ws1.Range("F8").Value = 1039
$ws1.Range("F9").Value = 401
$ws1.Range("F11").Value = 135
$ws1.Range("F12").Value = 11881
$ws1.Range("F13").Value = 36
$ws1.Range("F14").Value = 50
$ws1.Range("F15").Value = 11560
$ws1.Range("F16").Value = 4711
$ws1.Range("F17").Value = 501
$ws1.Range("F18").Value = 58
$ws1.Range("F20").Value = 83
$ws1.Range("F21").Value = 919
$ws1.Range("F24").Value = 46

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1252
$ws4.Range("F3").Value = 17108
$ws4.Range("F5").Value = 1665
$ws4.Range("F8").Value = 1039
$ws4.Range("F9").Value = 401
$ws4.Range("F11").Value = 135
$ws4.Range("F14").Value = 11881
$ws4.Range("F15").Value = 36
$ws4.Range("F16").Value = 50
$ws4.Range("F17").Value = 11560
$ws4.Range("F18").Value = 4711
$ws4.Range("F19").Value = 501
$ws4.Range("F20").Value = 58
$ws4.Range("F22").Value = 83
$ws4.Range("F23").Value = 919
$ws4.Range("F26").Value = 46
